$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 4 ("Request Type:") by copying the formatting of
# row 3 ("Requestor Signature:") which has the identical layout, then
# re-label it. This shifts every row from the old row 4 onward down by one,
# matching the rest of the diff (merged cells, dimension, print area, etc.)
$ws.Rows("3").Copy()
$ws.Rows("4").Insert()

# The merge state of the copied row doesn't come across automatically, so
# reapply it explicitly on the newly inserted row.
$ws.Range("A4:B4").Merge()
$ws.Range("C4:H4").Merge()

# Match row 3's height (27pt, custom height) and set the new label text.
$ws.Rows("4").RowHeight = $ws.Rows("3").RowHeight
$ws.Range("A4").Value = "Request Type:"

# The printed area grew by one row because of the inserted row.
$ws.PageSetup.PrintArea = '$A$1:$I$33'

# Match the final selection state recorded in the saved file.
$ws.Range("L1").Select()
